# Helper: assign a horizontal (1 x N) block of values to a range in one shot,
# building the 2-D SAFEARRAY Excel's Range.Value setter expects.
function Set-RowRange {
    param($ws, $rangeAddr, $values)
    $n = $values.Length
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $ws.Range($rangeAddr).Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gains a new trailing observation (K2).
$ws.Range("K2").Value = 0.0652635207166864

# Rows 3-14: each row's series shifts one column to the left (the oldest
# value in B is dropped) and a new value is appended in K.
Set-RowRange $ws "B3:K3" @(7.869792681105293,-10.45921331889471,-2.319131318894706,-0.5308223188947059,-3.708941318894706,-4.004270318894706,-1.489825318894706,-2.419547318894706,-2.184539318894706,-2.491675318894706)
Set-RowRange $ws "B4:K4" @(-17.928363649752,-9.788281649751999,-7.999972649751999,-11.178091649752,-11.473420649752,-8.958975649751999,-9.888697649751998,-9.653689649752,-9.960825649752,-9.450149649751999)
Set-RowRange $ws "B5:K5" @(8.382980150385469,10.17128915038547,6.99317015038547,6.69784115038547,9.21228615038547,8.28256415038547,8.517572150385469,8.210436150385469,8.721112150385469,8.433745150385469)
Set-RowRange $ws "B6:K6" @(2.03161760340632,-1.146501396593679,-1.441830396593679,1.072614603406321,0.1428926034063205,0.3779006034063205,0.07076460340632051,0.5814406034063205,0.2940736034063205,0.3220726034063205)
Set-RowRange $ws "B7:K7" @(-3.294071868393072,-3.589400868393072,-1.074955868393072,-2.004677868393072,-1.769669868393072,-2.076805868393072,-1.566129868393072,-1.853496868393072,-1.825497868393072,-1.732849868393072)
Set-RowRange $ws "B8:K8" @(-0.1493708818697419,2.365074118130258,1.435352118130258,1.670360118130258,1.363224118130258,1.873900118130258,1.586533118130258,1.614532118130258,1.707180118130258,1.920156118130258)
Set-RowRange $ws "B9:K9" @(2.426248296544459,1.496526296544459,1.731534296544459,1.424398296544459,1.935074296544459,1.647707296544459,1.675706296544459,1.768354296544459,1.981330296544459,1.480670296544459)
Set-RowRange $ws "B10:K10" @(-0.5449296406117954,-0.3099216406117954,-0.6170576406117954,-0.1063816406117954,-0.3937486406117954,-0.3657496406117954,-0.2731016406117954,-0.06012564061179543,-0.5607856406117955,-0.2710286406117954)
Set-RowRange $ws "B11:K11" @(0.2726728964218456,-0.03446310357815441,0.4762128964218456,0.1888458964218456,0.2168448964218456,0.3094928964218456,0.5224688964218456,0.02180889642184558,0.3115658964218456,0.05367389642184558)
Set-RowRange $ws "B12:K12" @(-0.1301381443875122,0.3805378556124878,0.09317085561248779,0.1211698556124878,0.2138178556124878,0.4267938556124878,-0.07386614438751221,0.2158908556124878,-0.04200114438751221,0.3403798556124878)
Set-RowRange $ws "B13:K13" @(0.6409568926112106,0.3535898926112106,0.3815888926112106,0.4742368926112106,0.6872128926112107,0.1865528926112106,0.4763098926112106,0.2184178926112106,0.6007988926112107,0.1895217986112106)
Set-RowRange $ws "B14:K14" @(-0.4818651035472806,-0.4538661035472806,-0.3612181035472806,-0.1482421035472806,-0.6489021035472806,-0.3591451035472806,-0.6170371035472806,-0.2346561035472806,-0.6459331975472806,-0.3589771035472806)

# Rows 15-24: same left-shift, but these rows shrink (no new trailing value,
# so the rightmost cell of each is cleared below).
Set-RowRange $ws "B15:J15" @(-0.3537867436446591,-0.2611387436446591,-0.04816274364465911,-0.5488227436446591,-0.2590657436446591,-0.5169577436446591,-0.1345767436446591,-0.5458538376446591,-0.2588977436446591)
Set-RowRange $ws "B16:I16" @(0.09264776243503714,0.3056237624350371,-0.1950362375649629,0.09472076243503715,-0.1631712375649629,0.2192097624350371,-0.1920673315649629,0.09488876243503713)
Set-RowRange $ws "B17:H17" @(0.1574463720025918,-0.3432136279974082,-0.05345662799740819,-0.3113486279974082,0.0710323720025918,-0.3402447219974082,-0.0532886279974082)
Set-RowRange $ws "B18:G18" @(-0.5006596170015631,-0.2109026170015632,-0.4687946170015632,-0.08641361700156319,-0.4976907110015631,-0.2107346170015632)
Set-RowRange $ws "B19:F19" @(0.2897568395245076,0.03186483952450761,0.4142458395245076,0.002968745524507627,0.2899248395245076)
Set-RowRange $ws "B20:E20" @(-0.2484930410109615,0.1338879589890384,-0.2773891350109615,0.009566958989038449)
Set-RowRange $ws "B21:D21" @(0.3058628168340501,-0.1054142771659499,0.1815418168340501)
Set-RowRange $ws "B22:C22" @(-0.4112768169122814,-0.1243207229122814)
$ws.Range("B23").Value = 0.4880092297750048

# Drop the now-unused trailing cell in each shrinking row.
$ws.Range("K15").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("H18").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("B24").ClearContents()
